# Auto-generated: reorders species-record rows 3-25 on the active sheet
# to match the target permutation (same 23 records, new row positions),
# including moving the M-column note cells with their owning record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{}
$rows[3] = @{
    A = 99564944
    B = 96334
    D = 'VU'
    E = 220787
    F = 'Knärot'
    G = 'Goodyera repens'
    H = '(L.) R. Br.'
    M = $null
    Q = 595280.9125479215
    R = 6912292.529634375
}
$rows[4] = @{
    A = 99564937
    B = 89392
    D = 'NT'
    E = 1202
    F = 'Ullticka'
    G = 'Phellinidium ferrugineofuscum'
    H = '(P.Karst.) Fiasson & Niemelä'
    M = $null
    Q = 595718.462900955
    R = 6912331.040771263
}
$rows[5] = @{
    A = 99564923
    B = 89376
    D = 'LC'
    E = 4660
    F = 'Rävticka'
    G = 'Inocutis rheades'
    H = '(Pers.) Fiasson & Niemelä'
    M = $null
    Q = 595468.7437919002
    R = 6912221.153764501
}
$rows[6] = @{
    A = 99564927
    B = 78503
    D = 'LC'
    E = 6456
    F = 'Skinnlav'
    G = 'Leptogium saturninum'
    H = '(Dicks.) Nyl.'
    M = $null
    Q = 595451.5082853716
    R = 6912384.801054608
}
$rows[7] = @{
    A = 99564931
    B = 78503
    D = 'LC'
    E = 6456
    F = 'Skinnlav'
    G = 'Leptogium saturninum'
    H = '(Dicks.) Nyl.'
    M = $null
    Q = 595438.9989911841
    R = 6912382.584794725
}
$rows[8] = @{
    A = 99564925
    B = 89392
    D = 'NT'
    E = 1202
    F = 'Ullticka'
    G = 'Phellinidium ferrugineofuscum'
    H = '(P.Karst.) Fiasson & Niemelä'
    M = $null
    Q = 595476.4437787337
    R = 6912245.08772236
}
$rows[9] = @{
    A = 99564928
    B = 78569
    D = 'NT'
    E = 6458
    F = 'Lunglav'
    G = 'Lobaria pulmonaria'
    H = '(L.) Hoffm.'
    M = $null
    Q = 595582.6011734826
    R = 6912474.092091525
}
$rows[10] = @{
    A = 99564932
    B = 77506
    D = 'NT'
    E = 6425
    F = 'Garnlav'
    G = 'Alectoria sarmentosa'
    H = '(Ach.) Ach.'
    M = $null
    Q = 595581.4055642756
    R = 6912483.357140777
}
$rows[11] = @{
    A = 99564930
    B = 78569
    D = 'NT'
    E = 6458
    F = 'Lunglav'
    G = 'Lobaria pulmonaria'
    H = '(L.) Hoffm.'
    M = $null
    Q = 595546.9585193637
    R = 6912287.086805391
}
$rows[12] = @{
    A = 99564943
    B = 77259
    D = 'NT'
    E = 228912
    F = 'Mörk kolflarnlav'
    G = 'Carbonicola myrmecina'
    H = '(Ach.) Bendiksby & Timdal'
    M = $null
    Q = 595659.09430371
    R = 6912339.109678851
}
$rows[13] = @{
    A = 99564940
    B = 78503
    D = 'LC'
    E = 6456
    F = 'Skinnlav'
    G = 'Leptogium saturninum'
    H = '(Dicks.) Nyl.'
    M = $null
    Q = 595470.5255192126
    R = 6912223.994488954
}
$rows[14] = @{
    A = 99564945
    B = 56395
    D = 'NT'
    E = 100109
    F = 'Tretåig hackspett'
    G = 'Picoides tridactylus'
    H = '(Linnaeus, 1758)'
    M = 'äldre spår'
    Q = 595778.4600258654
    R = 6912463.877241801
}
$rows[15] = @{
    A = 99564936
    B = 78569
    D = 'NT'
    E = 6458
    F = 'Lunglav'
    G = 'Lobaria pulmonaria'
    H = '(L.) Hoffm.'
    M = $null
    Q = 595471.9213900227
    R = 6912224.034264626
}
$rows[16] = @{
    A = 99564934
    B = 56395
    D = 'NT'
    E = 100109
    F = 'Tretåig hackspett'
    G = 'Picoides tridactylus'
    H = '(Linnaeus, 1758)'
    M = 'födosökande'
    Q = 595782.0215915864
    R = 6912420.737332884
}
$rows[17] = @{
    A = 99564924
    B = 78527
    D = 'LC'
    E = 229497
    F = 'Korallblylav'
    G = 'Parmeliella triptophylla'
    H = '(Ach.) Müll.Arg.'
    M = $null
    Q = 595505.2006268308
    R = 6912265.436834001
}
$rows[18] = @{
    A = 99564926
    B = 78569
    D = 'NT'
    E = 6458
    F = 'Lunglav'
    G = 'Lobaria pulmonaria'
    H = '(L.) Hoffm.'
    M = $null
    Q = 595657.1934100311
    R = 6912340.450337817
}
$rows[19] = @{
    A = 99564933
    B = 89356
    D = 'LC'
    E = 5447
    F = 'Vedticka'
    G = 'Fuscoporia viticola'
    H = '(Schwein.) Murrill'
    M = $null
    Q = 595465.3450126103
    R = 6912258.721109796
}
$rows[20] = @{
    A = 99564942
    B = 78602
    D = 'LC'
    E = 6463
    F = 'Bårdlav'
    G = 'Nephroma parile'
    H = '(Ach.) Ach.'
    M = $null
    Q = 595451.5744691773
    R = 6912382.478090141
}
$rows[21] = @{
    A = 99564929
    B = 78596
    D = 'LC'
    E = 6462
    F = 'Stuplav'
    G = 'Nephroma bellum'
    H = '(Spreng.) Tuck.'
    M = $null
    Q = 595509.0950829939
    R = 6912259.503016504
}
$rows[22] = @{
    A = 99564938
    B = 78569
    D = 'NT'
    E = 6458
    F = 'Lunglav'
    G = 'Lobaria pulmonaria'
    H = '(L.) Hoffm.'
    M = $null
    Q = 595794.0918679656
    R = 6912422.012381156
}
$rows[23] = @{
    A = 99564939
    B = 89392
    D = 'NT'
    E = 1202
    F = 'Ullticka'
    G = 'Phellinidium ferrugineofuscum'
    H = '(P.Karst.) Fiasson & Niemelä'
    M = $null
    Q = 595420.0288524196
    R = 6912378.324677907
}
$rows[24] = @{
    A = 99564941
    B = 89673
    D = 'NT'
    E = 658
    F = 'Rosenticka'
    G = 'Rhodofomes roseus'
    H = '(Alb. & Schwein.) Kotl. & Pouzar'
    M = $null
    Q = 595711.9895855145
    R = 6912345.735149743
}
$rows[25] = @{
    A = 99564935
    B = 89832
    D = 'VU'
    E = 1209
    F = 'Rynkskinn'
    G = 'Phlebia centrifuga'
    H = 'P.Karst.'
    M = $null
    Q = 595713.3854063898
    R = 6912345.775026423
}

foreach ($r in $rows.Keys) {
    $rec = $rows[$r]
    $ws.Range("A$r").Value = $rec.A
    $ws.Range("B$r").Value = $rec.B
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("E$r").Value = $rec.E
    $ws.Range("F$r").Value = $rec.F
    $ws.Range("G$r").Value = $rec.G
    $ws.Range("H$r").Value = $rec.H
    if ($null -eq $rec.M) {
        $ws.Range("M$r").Value = ""
    } else {
        $ws.Range("M$r").Value = $rec.M
    }
    $ws.Range("Q$r").Value = $rec.Q
    $ws.Range("R$r").Value = $rec.R
}
